# Apply updated inclusion/nucleus counts and recalculated ratios
# (analyze 120524 and improve algorithm)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 1K_PFF_ADAMTS19_01
$ws.Range("B2").Value = 29
$ws.Range("D2").Value = 2.9

# Row 3: 1K_PFF_ADAMTS19_02
$ws.Range("B3").Value = 18
$ws.Range("D3").Value = 1.058823529411765

# Row 6: 1K_PFF_ADAMTS19_05
$ws.Range("B6").Value = 82
$ws.Range("D6").Value = 41

# Row 7: 1K_PFF_ADAMTS19_06
$ws.Range("B7").Value = 20
$ws.Range("D7").Value = 1.428571428571429

# Row 8: 1K_PFF_ADAMTS19_07
$ws.Range("B8").Value = 27
$ws.Range("D8").Value = 2.076923076923077

# Row 9: 1K_PFF_ADAMTS19_08
$ws.Range("B9").Value = 4
$ws.Range("D9").Value = 0.5714285714285714

# Row 10: 1K_PFF_ADAMTS19_09
$ws.Range("B10").Value = 0
$ws.Range("D10").Value = 0

# Row 12: 1K_PFF_SCR_01
$ws.Range("B12").Value = 106
$ws.Range("D12").Value = 17.66666666666667

# Row 13: 1K_PFF_SCR_02
$ws.Range("B13").Value = 9
$ws.Range("D13").Value = 0.3333333333333333

# Row 15: 1K_PFF_SCR_04
$ws.Range("B15").Value = 1
$ws.Range("D15").Value = 0.03703703703703703

# Row 23: 1K_PFF_TAX1BP1_02
$ws.Range("B23").Value = 24
$ws.Range("D23").Value = 2.181818181818182
